# contratos-4-2015.xlsx formatting fix
# ------------------------------------
# 1) A handful of "Razon social"/"Nombre Fantasia" entries used a comma
#    as a separator between multiple contractor names; normalize those
#    separators to a period (and drop stray dots in "S.H." -> "SH").
# 2) The "Importe" column (H) was scraped with Spanish/Latin-style
#    thousands/decimal separators ("1.234,56") but stored as plain text.
#    Rewrite every value as a plain decimal string ("1234.56").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Razon social / Nombre Fantasia separator fixes ---
$nameCells = @(
    "E80",
    "E83",
    "E121",
    "E180",
    "E195",
    "F164"
)
$nameValues = @(
    "FERNANDEZ. MARIO HUGO",
    "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO",
    "RICCOTTI. MARIANA EDITH",
    "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH",
    "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH",
    "MERCANZINI. GASTON ARIEL"
)
for ($i = 0; $i -lt $nameCells.Length; $i++) {
    $ws.Range($nameCells[$i]).Value = $nameValues[$i]
}

# --- 2) "Importe" column (H2:H255): re-write Spanish-formatted numeric
#    text as plain decimal text. The column holds TEXT (not numbers), so
#    the cells are switched to the Text format before writing - otherwise
#    Excel would parse a string like "6035.00" back into the number 6035
#    and drop the trailing ".00" - then switched back to the workbook's
#    default ("Normal") style so no stray formatting is left behind.
$importeCells = @(
    "H2",
    "H3",
    "H4",
    "H5",
    "H6",
    "H7",
    "H8",
    "H9",
    "H10",
    "H11",
    "H12",
    "H13",
    "H14",
    "H15",
    "H16",
    "H17",
    "H18",
    "H19",
    "H20",
    "H21",
    "H22",
    "H23",
    "H24",
    "H25",
    "H26",
    "H27",
    "H28",
    "H29",
    "H30",
    "H31",
    "H32",
    "H33",
    "H34",
    "H35",
    "H36",
    "H37",
    "H38",
    "H39",
    "H40",
    "H41",
    "H42",
    "H43",
    "H44",
    "H45",
    "H46",
    "H47",
    "H48",
    "H49",
    "H50",
    "H51",
    "H52",
    "H53",
    "H54",
    "H55",
    "H56",
    "H57",
    "H58",
    "H59",
    "H60",
    "H61",
    "H62",
    "H63",
    "H64",
    "H65",
    "H66",
    "H67",
    "H68",
    "H69",
    "H70",
    "H71",
    "H72",
    "H73",
    "H74",
    "H75",
    "H76",
    "H77",
    "H78",
    "H79",
    "H80",
    "H81",
    "H82",
    "H83",
    "H84",
    "H85",
    "H86",
    "H87",
    "H88",
    "H89",
    "H90",
    "H91",
    "H92",
    "H93",
    "H94",
    "H95",
    "H96",
    "H97",
    "H98",
    "H99",
    "H100",
    "H101",
    "H102",
    "H103",
    "H104",
    "H105",
    "H106",
    "H107",
    "H108",
    "H109",
    "H110",
    "H111",
    "H112",
    "H113",
    "H114",
    "H115",
    "H116",
    "H117",
    "H118",
    "H119",
    "H120",
    "H121",
    "H122",
    "H123",
    "H124",
    "H125",
    "H126",
    "H127",
    "H128",
    "H129",
    "H130",
    "H131",
    "H132",
    "H133",
    "H134",
    "H135",
    "H136",
    "H137",
    "H138",
    "H139",
    "H140",
    "H141",
    "H142",
    "H143",
    "H144",
    "H145",
    "H146",
    "H147",
    "H148",
    "H149",
    "H150",
    "H151",
    "H152",
    "H153",
    "H154",
    "H155",
    "H156",
    "H157",
    "H158",
    "H159",
    "H160",
    "H161",
    "H162",
    "H163",
    "H164",
    "H165",
    "H166",
    "H167",
    "H168",
    "H169",
    "H170",
    "H171",
    "H172",
    "H173",
    "H174",
    "H175",
    "H176",
    "H177",
    "H178",
    "H179",
    "H180",
    "H181",
    "H182",
    "H183",
    "H184",
    "H185",
    "H186",
    "H187",
    "H188",
    "H189",
    "H190",
    "H191",
    "H192",
    "H193",
    "H194",
    "H195",
    "H196",
    "H197",
    "H198",
    "H199",
    "H200",
    "H201",
    "H202",
    "H203",
    "H204",
    "H205",
    "H206",
    "H207",
    "H208",
    "H209",
    "H210",
    "H211",
    "H212",
    "H213",
    "H214",
    "H215",
    "H216",
    "H217",
    "H218",
    "H219",
    "H220",
    "H221",
    "H222",
    "H223",
    "H224",
    "H225",
    "H226",
    "H227",
    "H228",
    "H229",
    "H230",
    "H231",
    "H232",
    "H233",
    "H234",
    "H235",
    "H236",
    "H237",
    "H238",
    "H239",
    "H240",
    "H241",
    "H242",
    "H243",
    "H244",
    "H245",
    "H246",
    "H247",
    "H248",
    "H249",
    "H250",
    "H251",
    "H252",
    "H253",
    "H254",
    "H255"
)
$importeValues = @(
    "6035.00",
    "175000.00",
    "180000.00",
    "1049.40",
    "421273.43",
    "17424.00",
    "365311.10",
    "4555.00",
    "3612.00",
    "3210.00",
    "433089.27",
    "83600.20",
    "105894.11",
    "3000.00",
    "8097.00",
    "1774.59",
    "13534.73",
    "14045.95",
    "3751.00",
    "26040.00",
    "13163.92",
    "25506.04",
    "4500.00",
    "1800.00",
    "4312.00",
    "6080.00",
    "262.54",
    "507.88",
    "9630.00",
    "72790.96",
    "884.57",
    "9699.00",
    "2806.00",
    "113632.36",
    "882.80",
    "3279493.80",
    "31514.00",
    "833.28",
    "5299.00",
    "2582.81",
    "45428.00",
    "51250.57",
    "2910.00",
    "894.80",
    "510.00",
    "171472.00",
    "47.68",
    "14341.66",
    "4257.64",
    "444.71",
    "252.50",
    "207583.37",
    "1521.65",
    "33296.05",
    "20830.98",
    "4466.03",
    "95.40",
    "2457.02",
    "979.80",
    "1508.00",
    "1312.71",
    "791.00",
    "3828.00",
    "1516.00",
    "1043.57",
    "1941.28",
    "2865.00",
    "8425.00",
    "4200.00",
    "79.00",
    "25815.00",
    "875.00",
    "3450.00",
    "60.00",
    "17890.00",
    "2148.00",
    "14000.00",
    "179.00",
    "23055.00",
    "34585.20",
    "1400.00",
    "8540.00",
    "100.00",
    "6000.00",
    "11280.00",
    "273830.58",
    "3999.99",
    "4600.00",
    "371897.61",
    "59551.45",
    "47.93",
    "181.85",
    "21294.95",
    "85.32",
    "777.80",
    "5474.50",
    "2607.98",
    "2400.00",
    "4637.90",
    "7665.00",
    "128.00",
    "400.00",
    "1709.00",
    "2008.00",
    "75.00",
    "3720.00",
    "250.00",
    "1800.00",
    "28341.90",
    "242.00",
    "80.00",
    "9835.72",
    "130.00",
    "41.00",
    "95.60",
    "200.00",
    "85.00",
    "540.00",
    "175.00",
    "2000.00",
    "360.00",
    "30000.00",
    "204999.75",
    "88240.00",
    "2375.00",
    "6260.00",
    "2500.00",
    "4360.00",
    "3980.90",
    "337590.00",
    "3435.00",
    "26278.00",
    "4000.00",
    "4000.00",
    "2385.00",
    "500.00",
    "32400.00",
    "538.00",
    "918.23",
    "87.36",
    "200.00",
    "1733.00",
    "400.00",
    "721.26",
    "8156.78",
    "242.50",
    "10000.35",
    "63.00",
    "29124.00",
    "11486.94",
    "2300.00",
    "1050.00",
    "1500.00",
    "55158.18",
    "3146.00",
    "1200.00",
    "1105.00",
    "1440.00",
    "2000.00",
    "800.00",
    "1000.00",
    "13500.00",
    "9000.00",
    "4000.00",
    "950.00",
    "1200.00",
    "10359.61",
    "1862.00",
    "3000.00",
    "600.00",
    "1260.00",
    "24080.00",
    "1300.00",
    "1300.00",
    "4600.00",
    "300.00",
    "2400.00",
    "9674.00",
    "850.00",
    "300.00",
    "2060.00",
    "322.00",
    "5110.00",
    "2000.00",
    "4080.00",
    "435.00",
    "912.84",
    "2718.80",
    "1480.00",
    "3332.96",
    "1656.00",
    "8593.80",
    "6640.00",
    "17640.00",
    "28.98",
    "1719.00",
    "1171.15",
    "300.00",
    "1251.80",
    "74125.00",
    "2350.00",
    "205827.00",
    "72015.00",
    "6146.49",
    "14506.30",
    "2554.00",
    "700.00",
    "368225.89",
    "3371.20",
    "5082.97",
    "42700.00",
    "1785.20",
    "777315.46",
    "148986.00",
    "600.00",
    "179336.00",
    "206612.22",
    "185000.00",
    "218700.00",
    "424303.73",
    "200000.00",
    "40414.70",
    "95677.00",
    "54719.60",
    "121781.25",
    "340221.00",
    "245952.00",
    "44991.01",
    "224748.00",
    "249553.80",
    "100000.00",
    "137948.00",
    "400569.60",
    "7820524.96",
    "18310.00",
    "507047.10",
    "78000.00",
    "11964.00",
    "58430.00",
    "1900.00",
    "6900.00",
    "108500.00",
    "350.00",
    "4600.00",
    "7450.00",
    "4000.00",
    "1750.00",
    "18225.75",
    "18000.00",
    "1783.54",
    "384900.00",
    "396.00",
    "2658.00",
    "4600.00"
)

$importeRange = $ws.Range("H2:H255")
$importeRange.NumberFormat = "@"
for ($i = 0; $i -lt $importeCells.Length; $i++) {
    $ws.Range($importeCells[$i]).Value = $importeValues[$i]
}
$importeRange.Style = "Normal"
